# FACET vectors generated for 339-377 videos
#
# The "Annotation-Supervised" sheet lists the interview question asked in
# each annotated row alongside its polarity/label. Row 21 had been sharing
# the exact same shared-string text as rows 19-20 ("how close are you to
# your family"), even though it is really a distinct annotation instance
# (a separate video in the 339-377 batch). Disambiguate it by giving it its
# own text value (a trailing-space variant), which forces Excel to create a
# new, separate shared-string entry instead of continuing to alias the
# shared one used by rows 19 and 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Annotation-Supervised")

$ws.Range("A21").Value = "how close are you to your family "
